# "Rephrase multiple teams to tech orga"
#
# On slides 2-20 there is a small two-line label shape (named "TextBox ...")
# nested one level inside a top-level group shape (named "Group ...").
# It originally read:
#     "Multiple "
#     "Teams"
# and needs to become:
#     "Tech "
#     "Orga"
# The shape uses <a:spAutoFit/>, so PowerPoint also re-flows its bounding
# box to the new (narrower) text; the target width is 504433 EMU (the
# height is unchanged at 461665 EMU).

function Set-TechOrgaLabel($shp) {
    $tr = $shp.TextFrame.TextRange

    # Replace paragraph 1 ("Multiple ") with "Tech ", inserting the new
    # text next to the still-present paragraph 2 so it inherits that
    # paragraph's run formatting (incl. the "dirty" bookkeeping flag),
    # same as PowerPoint does when you type over existing text.
    $para1 = $tr.Paragraphs(1, 1)
    $para1.Delete()
    $tr.InsertBefore("Tech `r")

    # Replace paragraph 2's run ("Teams") with "Orga", keeping the
    # paragraph mark (and its formatting) in place.
    $para2 = $tr.Paragraphs(2, 1)
    $run2 = $para2.Runs(1, 1)
    $run2.Delete()
    $tr.InsertAfter("Orga")

    # The shape auto-fits to its (now shorter) text; pin the resulting
    # size to match what PowerPoint computed (745717 x 461665 -> 504433 x
    # 461665 EMU). Left/Top are untouched.
    $shp.Width = 39.71914
    $shp.Height = 36.3516
}

$p = $ppt.ActivePresentation

$targets = @(
    @{slide = 2;  grp = "Group 87";  tb = "TextBox 70"},
    @{slide = 3;  grp = "Group 7";   tb = "TextBox 199"},
    @{slide = 4;  grp = "Group 11";  tb = "TextBox 244"},
    @{slide = 5;  grp = "Group 5";   tb = "TextBox 249"},
    @{slide = 6;  grp = "Group 10";  tb = "TextBox 206"},
    @{slide = 7;  grp = "Group 5";   tb = "TextBox 161"},
    @{slide = 8;  grp = "Group 6";   tb = "TextBox 161"},
    @{slide = 9;  grp = "Group 6";   tb = "TextBox 205"},
    @{slide = 10; grp = "Group 8";   tb = "TextBox 161"},
    @{slide = 11; grp = "Group 8";   tb = "TextBox 161"},
    @{slide = 12; grp = "Group 12";  tb = "TextBox 206"},
    @{slide = 13; grp = "Group 138"; tb = "TextBox 154"},
    @{slide = 14; grp = "Group 7";   tb = "TextBox 161"},
    @{slide = 15; grp = "Group 8";   tb = "TextBox 160"},
    @{slide = 16; grp = "Group 9";   tb = "TextBox 161"},
    @{slide = 17; grp = "Group 7";   tb = "TextBox 161"},
    @{slide = 18; grp = "Group 7";   tb = "TextBox 161"},
    @{slide = 19; grp = "Group 7";   tb = "TextBox 161"},
    @{slide = 20; grp = "Group 7";   tb = "TextBox 161"}
)

foreach ($row in $targets) {
    $slide = $p.Slides.Item($row.slide)
    $grp = $slide.Shapes.Item($row.grp)
    $shp = $grp.GroupItems.Item($row.tb)
    Set-TechOrgaLabel $shp
}
